# Turkey 1 Lig - base update (17-02-2024 22:47)
# The edit re-sorts a handful of adjacent (or near-adjacent) match rows by
# swapping their full data payload (columns B:AC) while keeping the row's
# running index in column A fixed, and removes the final obsolete row
# (old row 354) by shifting the last few rows up by one and deleting the
# now-duplicated tail row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valA = $rangeA.Value()
    $valB = $rangeB.Value()
    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# Pairs of rows whose full records (everything except the col-A running
# index) need to be swapped.
Swap-Rows 24 25
Swap-Rows 54 55
Swap-Rows 74 75
Swap-Rows 125 127
Swap-Rows 179 180
Swap-Rows 195 196
Swap-Rows 208 209
Swap-Rows 214 215
Swap-Rows 219 220
Swap-Rows 250 251
Swap-Rows 269 270
Swap-Rows 288 289
Swap-Rows 304 305
Swap-Rows 335 336

# Rows 349-353 each shift up to take the next row's record, and the
# now-redundant last row (354) is removed entirely.
for ($r = 349; $r -le 353; $r++) {
    $next = $r + 1
    $src = $ws.Range("B$next`:AC$next")
    $dst = $ws.Range("B$r`:AC$r")
    $dst.Value = $src.Value()
}

$ws.Rows(354).Delete()
